# Update projections script to read in round 2 UI override
# Updates columns C (impact), E (total), H (consumption) for rows 86-94 on Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    86 = @{ C = 3.738774692380103;  E = -3.026985394548855; H = -2.860268089315993 }
    87 = @{ C = 5.55998151019559;   E = 7.3698949329205;     H = 5.223186118352026 }
    88 = @{ C = 1.137011959877302;  E = -3.402100616507544; H = -2.963613217270506 }
    89 = @{ C = -1.206409713280085; E = -5.766447774984444; H = -4.738360023252155 }
    90 = @{ C = -0.9480192944861663;E = -1.993423719373181; H = -1.857206675026766 }
    91 = @{ C = -3.344155676585455; E = -2.214650595476656; H = -1.607805366189948 }
    92 = @{ C = -3.009298306419569; E = -2.062671135843999; H = -1.993773408565533 }
    93 = @{ C = -1.740247042006141; E = -0.6902427173307312;H = -0.5832714240855903 }
    94 = @{ C = -1.287716044876893; E = -0.1832997308561881;H = 0.009126076209881323 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("H$row").Value = $vals.H
}
